# ============================================================
# Commit: "added norm100 and befor market opening condition"
# - Rename header AA1: stop_limit_sell_order_id -> stop_limit_order_id
# - Append 7 new trade rows (270-276), extending dimension to A1:AA276
# ============================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header cell AA1 ---
$ws.Range("AA1").Value = 'stop_limit_order_id'

# ---------------- Row 270 ----------------
# Clone cell formatting (format-only paste) for the styled columns from
# the last existing data row, so new cells reuse the same style indexes
# instead of minting new ones.
$ws.Range("A269").Copy()
$ws.Range("A270").PasteSpecial(-4122)
$ws.Range("D269").Copy()
$ws.Range("D270").PasteSpecial(-4122)
$ws.Range("D269").Copy()
$ws.Range("H270").PasteSpecial(-4122)

# Set cell values
$ws.Range("A270").Value = 268
$ws.Range("B270").Value = 290
$ws.Range("C270").Value = 'FTNT'
$ws.Range("D270").Value = 45580.18006612269
$ws.Range("E270").Value = 82.5492
$ws.Range("F270").Value = 3136.8696
$ws.Range("G270").Value = 1.1
$ws.Range("H270").Value = 45580.28755964121
$ws.Range("I270").Value = 82.66030000000001
$ws.Range("J270").Value = 3141.0914
$ws.Range("K270").Value = 1.2
$ws.Range("L270").Value = 38
$ws.Range("M270").Value = 'completed'
$ws.Range("N270").Value = 1.005
$ws.Range("O270").Value = 0.995
$ws.Range("P270").Value = 1.006
$ws.Range("Q270").Value = 0.02
$ws.Range("R270").Value = 1.921800000000257
$ws.Range("S270").Value = 'FA199DDF65447B2000'
$ws.Range("W270").Value = 'FA199DE6C57A44A000'
$ws.Range("X270").Value = '1230'
$ws.Range("Y270").Value = 'mv :-1.14, mv_2m:-0.15,      mv_5m : -0.21, mv_30m : -0.60, mv_60m: -0.05'
$ws.Range("Z270").Value = 'Сиднейское время (лето)'

# Materialize empty-but-present cells (no value, default style)
$ws.Range("T269").Copy()
$ws.Range("T270").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("U270").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("V270").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("AA270").PasteSpecial(-4122)

# ---------------- Row 271 ----------------
# Clone cell formatting (format-only paste) for the styled columns from
# the last existing data row, so new cells reuse the same style indexes
# instead of minting new ones.
$ws.Range("A269").Copy()
$ws.Range("A271").PasteSpecial(-4122)
$ws.Range("D269").Copy()
$ws.Range("D271").PasteSpecial(-4122)
$ws.Range("D269").Copy()
$ws.Range("H271").PasteSpecial(-4122)

# Set cell values
$ws.Range("A271").Value = 269
$ws.Range("B271").Value = 293
$ws.Range("C271").Value = 'CMG'
$ws.Range("D271").Value = 45581.15869940972
$ws.Range("E271").Value = 59.72
$ws.Range("F271").Value = 3165.16
$ws.Range("G271").Value = 1.15
$ws.Range("H271").Value = 45581.28762678241
$ws.Range("I271").Value = 59.9001
$ws.Range("J271").Value = 3174.7053
$ws.Range("K271").Value = 1.101
$ws.Range("L271").Value = 53
$ws.Range("M271").Value = 'completed'
$ws.Range("N271").Value = 1.005
$ws.Range("O271").Value = 0.995
$ws.Range("P271").Value = 1.006
$ws.Range("Q271").Value = 0.02
$ws.Range("R271").Value = 7.294300000000224
$ws.Range("S271").Value = 'FA199F21F15A44A000'
$ws.Range("W271").Value = 'FA199F2BC9BBBB2000'
$ws.Range("X271").Value = '1230'
$ws.Range("Y271").Value = 'mv :-22.67, mv_2m:-0.50,      mv_5m : -1.72, mv_30m : 1.23, mv_60m: 3.43'
$ws.Range("Z271").Value = 'Сиднейское время (лето)'

# Materialize empty-but-present cells (no value, default style)
$ws.Range("T269").Copy()
$ws.Range("T271").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("U271").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("V271").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("AA271").PasteSpecial(-4122)

# ---------------- Row 272 ----------------
# Clone cell formatting (format-only paste) for the styled columns from
# the last existing data row, so new cells reuse the same style indexes
# instead of minting new ones.
$ws.Range("A269").Copy()
$ws.Range("A272").PasteSpecial(-4122)
$ws.Range("D269").Copy()
$ws.Range("D272").PasteSpecial(-4122)
$ws.Range("D269").Copy()
$ws.Range("H272").PasteSpecial(-4122)

# Set cell values
$ws.Range("A272").Value = 270
$ws.Range("B272").Value = 295
$ws.Range("C272").Value = 'BAC'
$ws.Range("D272").Value = 45582.17252115741
$ws.Range("E272").Value = 42.915
$ws.Range("F272").Value = 3175.71
$ws.Range("G272").Value = 1.21
$ws.Range("H272").Value = 45582.24687973379
$ws.Range("I272").Value = 42.975
$ws.Range("J272").Value = 3180.15
$ws.Range("K272").Value = 1.31
$ws.Range("L272").Value = 74
$ws.Range("M272").Value = 'completed'
$ws.Range("N272").Value = 1.005
$ws.Range("O272").Value = 0.995
$ws.Range("P272").Value = 1.006
$ws.Range("Q272").Value = 0.3
$ws.Range("R272").Value = 1.920000000000055
$ws.Range("S272").Value = 'FA19A070169E04A000'
$ws.Range("W272").Value = 'FA19A07897EE7B2000'
$ws.Range("X272").Value = '1230'
$ws.Range("Y272").Value = 'mv :-1.77, mv_2m:-0.15,      mv_5m : -0.28, mv_30m : -1.13, mv_60m: -2.30'
$ws.Range("Z272").Value = 'Сиднейское время (лето)'

# Materialize empty-but-present cells (no value, default style)
$ws.Range("T269").Copy()
$ws.Range("T272").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("U272").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("V272").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("AA272").PasteSpecial(-4122)

# ---------------- Row 273 ----------------
# Clone cell formatting (format-only paste) for the styled columns from
# the last existing data row, so new cells reuse the same style indexes
# instead of minting new ones.
$ws.Range("A269").Copy()
$ws.Range("A273").PasteSpecial(-4122)
$ws.Range("D269").Copy()
$ws.Range("D273").PasteSpecial(-4122)
$ws.Range("D269").Copy()
$ws.Range("H273").PasteSpecial(-4122)

# Set cell values
$ws.Range("A273").Value = 271
$ws.Range("B273").Value = 298
$ws.Range("C273").Value = 'ADSK'
$ws.Range("D273").Value = 45582.2124253125
$ws.Range("E273").Value = 287.83
$ws.Range("F273").Value = 3166.13
$ws.Range("G273").Value = 1.02
$ws.Range("H273").Value = 45582.28688711805
$ws.Range("I273").Value = 288.39
$ws.Range("J273").Value = 3172.29
$ws.Range("K273").Value = 1.12
$ws.Range("L273").Value = 11
$ws.Range("M273").Value = 'completed'
$ws.Range("N273").Value = 1.005
$ws.Range("O273").Value = 0.995
$ws.Range("P273").Value = 1.006
$ws.Range("Q273").Value = 0.02
$ws.Range("R273").Value = 4.02000000000031
$ws.Range("S273").Value = 'FA19A07D3D8B3B2000'
$ws.Range("W273").Value = 'FA19A081561684A000'
$ws.Range("X273").Value = '1230'
$ws.Range("Y273").Value = 'mv :-1.50, mv_2m:0.41,      mv_5m : 0.26, mv_30m : 0.60, mv_60m: 0.35'
$ws.Range("Z273").Value = 'Сиднейское время (лето)'

# Materialize empty-but-present cells (no value, default style)
$ws.Range("T269").Copy()
$ws.Range("T273").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("U273").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("V273").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("AA273").PasteSpecial(-4122)

# ---------------- Row 274 ----------------
# Clone cell formatting (format-only paste) for the styled columns from
# the last existing data row, so new cells reuse the same style indexes
# instead of minting new ones.
$ws.Range("A269").Copy()
$ws.Range("A274").PasteSpecial(-4122)
$ws.Range("D269").Copy()
$ws.Range("D274").PasteSpecial(-4122)
$ws.Range("D269").Copy()
$ws.Range("H274").PasteSpecial(-4122)

# Set cell values
$ws.Range("A274").Value = 274
$ws.Range("B274").Value = 301
$ws.Range("C274").Value = 'TT'
$ws.Range("D274").Value = 45583.1528390625
$ws.Range("E274").Value = 401.62
$ws.Range("F274").Value = 3212.96
$ws.Range("G274").Value = 1.01
$ws.Range("H274").Value = 45583.28773255787
$ws.Range("I274").Value = 401.35
$ws.Range("J274").Value = 3210.8
$ws.Range("K274").Value = 1.11
$ws.Range("L274").Value = 8
$ws.Range("M274").Value = 'completed'
$ws.Range("N274").Value = 1.005
$ws.Range("O274").Value = 0.995
$ws.Range("P274").Value = 1.006
$ws.Range("Q274").Value = 0.02
$ws.Range("R274").Value = -4.279999999999855
$ws.Range("S274").Value = 'FA19A1B3310AC4A000'
$ws.Range("W274").Value = 'FA19A1D24FB7BB2000'
$ws.Range("X274").Value = '1230'
$ws.Range("Y274").Value = 'mv :0.89, mv_2m:0.81,      mv_5m : 1.66, mv_30m : -1.27, mv_60m: -1.63'
$ws.Range("Z274").Value = 'Сиднейское время (лето)'

# Materialize empty-but-present cells (no value, default style)
$ws.Range("T269").Copy()
$ws.Range("T274").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("U274").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("V274").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("AA274").PasteSpecial(-4122)

# ---------------- Row 275 ----------------
# Clone cell formatting (format-only paste) for the styled columns from
# the last existing data row, so new cells reuse the same style indexes
# instead of minting new ones.
$ws.Range("A269").Copy()
$ws.Range("A275").PasteSpecial(-4122)
$ws.Range("D269").Copy()
$ws.Range("D275").PasteSpecial(-4122)
$ws.Range("D269").Copy()
$ws.Range("H275").PasteSpecial(-4122)

# Set cell values
$ws.Range("A275").Value = 275
$ws.Range("B275").Value = 302
$ws.Range("C275").Value = 'CHTR'
$ws.Range("D275").Value = 45583.16485706018
$ws.Range("E275").Value = 328.87
$ws.Range("F275").Value = 2959.83
$ws.Range("G275").Value = 1.02
$ws.Range("H275").Value = 25934
$ws.Range("J275").Value = 0
$ws.Range("K275").Value = 0
$ws.Range("L275").Value = 9
$ws.Range("M275").Value = 'bought'
$ws.Range("N275").Value = 1.005
$ws.Range("O275").Value = 0.995
$ws.Range("P275").Value = 1.006
$ws.Range("Q275").Value = 0.45
$ws.Range("R275").Value = 0
$ws.Range("S275").Value = 'FA19A1B7270C7B2000'
$ws.Range("X275").Value = '1230'
$ws.Range("Y275").Value = 'mv :3.49, mv_2m:0.07,      mv_5m : 0.03, mv_30m : 0.88, mv_60m: 0.35'
$ws.Range("Z275").Value = 'Сиднейское время (лето)'

# Materialize empty-but-present cells (no value, default style)
$ws.Range("T269").Copy()
$ws.Range("I275").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("T275").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("U275").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("V275").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("W275").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("AA275").PasteSpecial(-4122)

# ---------------- Row 276 ----------------
# Clone cell formatting (format-only paste) for the styled columns from
# the last existing data row, so new cells reuse the same style indexes
# instead of minting new ones.
$ws.Range("A269").Copy()
$ws.Range("A276").PasteSpecial(-4122)
$ws.Range("D269").Copy()
$ws.Range("D276").PasteSpecial(-4122)
$ws.Range("D269").Copy()
$ws.Range("H276").PasteSpecial(-4122)

# Set cell values
$ws.Range("A276").Value = 276
$ws.Range("B276").Value = 303
$ws.Range("C276").Value = 'AMAT'
$ws.Range("D276").Value = 45583.27386659722
$ws.Range("E276").Value = 203.21
$ws.Range("F276").Value = 3251.36
$ws.Range("G276").Value = 1.111
$ws.Range("H276").Value = 25934
$ws.Range("J276").Value = 0
$ws.Range("K276").Value = 0
$ws.Range("L276").Value = 16
$ws.Range("M276").Value = 'placed'
$ws.Range("N276").Value = 1.05
$ws.Range("O276").Value = 0.995
$ws.Range("P276").Value = 1.006
$ws.Range("Q276").Value = 0.45
$ws.Range("R276").Value = 0
$ws.Range("S276").Value = 'FA19A1DB14AABB2000'
$ws.Range("X276").Value = 'speed_norm100'
$ws.Range("Y276").Value = 'mv :0.16, mv_2m:-0.29,      mv_5m : -0.02, mv_30m : -1.16, mv_60m: -0.20'
$ws.Range("Z276").Value = 'Сиднейское время (лето)'

# Materialize empty-but-present cells (no value, default style)
$ws.Range("T269").Copy()
$ws.Range("I276").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("T276").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("U276").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("V276").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("W276").PasteSpecial(-4122)
$ws.Range("T269").Copy()
$ws.Range("AA276").PasteSpecial(-4122)
